$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Header rename (A1:D1)
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2. Capitalize connector words ("de", "del", "la", "las", "los", "el", "y") in place names
$nameFixes = @{
    "B6" = "Rincón De Romos"
    "B25" = "Amatenango De La Frontera"
    "B31" = "Comitán De Domínguez"
    "B37" = "Mazapa De Madero"
    "B39" = "Montecristo De Guerrero"
    "B45" = "Salto De Agua"
    "B46" = "San Cristóbal De Las Casas"
    "B67" = "Hidalgo Del Parral"
    "B86" = "San Juan De Sabinas"
    "B97" = "Villa De Álvarez"
    "A99" = "Ciudad De México"
    "B103" = "Cuajimalpa De Morelos"
    "A136" = "Estado De México"
    "B136" = "Acambay De Ruíz Castañeda"
    "B138" = "Almoloya De Alquisiras"
    "B139" = "Almoloya De Juárez"
    "B143" = "Atizapán De Zaragoza"
    "B150" = "Ecatepec De Morelos"
    "B155" = "Ixtapan De La Sal"
    "B156" = "Ixtapan Del Oro"
    "B162" = "Naucalpan De Juárez"
    "B166" = "San Felipe Del Progreso"
    "B167" = "Soyaniquilpan De Juárez"
    "B174" = "Tenango Del Valle"
    "B178" = "Tlalnepantla De Baz"
    "B182" = "Villa De Allende"
    "B191" = "Apaseo El Alto"
    "B192" = "Apaseo El Grande"
    "B199" = "Dolores Hidalgo Cuna De La Independencia Nacional"
    "B203" = "Jaral Del Progreso"
    "B210" = "Purísima Del Rincón"
    "B214" = "San Diego De La Unión"
    "B216" = "San Francisco Del Rincón"
    "B217" = "San Luis De La Paz"
    "B218" = "Santa Cruz De Juventino Rosas"
    "B220" = "Silao De La Victoria"
    "B224" = "Valle De Santiago"
    "B229" = "Acapulco De Juárez"
    "B230" = "Ajuchitlán Del Progreso"
    "B232" = "Atoyac De Álvarez"
    "B234" = "Chilpancingo De Los Bravo"
    "B238" = "Cutzamala De Pinzón"
    "B241" = "Huitzuco De Los Figueroa"
    "B242" = "Iguala De La Independencia"
    "B244" = "La Unión De Isidoro Montes De Oca"
    "B256" = "Taxco De Alarcón"
    "B258" = "Técpan De Galeana"
    "B263" = "Tlapa De Comonfort"
    "B272" = "Agua Blanca De Iturbide"
    "B276" = "Atotonilco El Grande"
    "B278" = "Cuautepec De Hinojosa"
    "B288" = "Pachuca De Soto"
    "B293" = "Tepeji Del Río De Ocampo"
    "B294" = "Tezontepec De Aldama"
    "B297" = "Tula De Allende"
    "B298" = "Tulancingo De Bravo"
    "B301" = "Zacualtipán De Ángeles"
    "B304" = "Acatlán De Juárez"
    "B305" = "Ahualulco De Mercado"
    "B308" = "Atotonilco El Alto"
    "B310" = "Autlán De Navarro"
    "B315" = "Cañadas De Obregón"
    "B323" = "Encarnación De Díaz"
    "B329" = "Huejuquilla El Alto"
    "B330" = "Ixtlahuacán Del Río"
    "B340" = "Ojuelos De Jalisco"
    "B344" = "San Cristóbal De La Barranca"
    "B347" = "San Martín De Bolaños"
    "B349" = "San Miguel El Alto"
    "B350" = "San Sebastián Del Oeste"
    "B353" = "Talpa De Allende"
    "B354" = "Tamazula De Gordiano"
    "B357" = "Teocuitatlán De Corona"
    "B358" = "Tepatitlán De Morelos"
    "B360" = "Tlajomulco De Zúñiga"
    "B366" = "Unión De San Antonio"
    "B367" = "Unión De Tula"
    "B368" = "Valle De Guadalupe"
    "B371" = "Yahualica De González Gallo"
    "B372" = "Zacoalco De Torres"
    "B375" = "Zapotitlán De Vadillo"
    "B376" = "Zapotlán El Grande"
    "B397" = "Coalcomán De Vázquez Pallares"
    "B399" = "Cojumatlán De Régules"
    "B481" = "Tlaltizapán De Zapata"
    "B489" = "Amatlán De Cañas"
    "B490" = "Bahía De Banderas"
    "B494" = "Ixtlán Del Río"
    "B501" = "Santa María Del Oro"
    "B512" = "San Nicolás De Los Garza"
    "B516" = "Coicoyán De Las Flores"
    "B517" = "Constancia Del Rosario"
    "B518" = "Heroica Ciudad De Huajuapan De León"
    "B519" = "Heroica Ciudad De Tlaxiaco"
    "B520" = "Heroica Ciudad De Juchitán De Zaragoza"
    "B522" = "Miahuatlán De Porfirio Díaz"
    "B523" = "Nejapa De Madero"
    "B524" = "Oaxaca De Juárez"
    "B525" = "Ocotlán De Morelos"
    "B572" = "Tataltepec De Valdés"
    "B573" = "Teotitlán De Flores Magón"
    "B574" = "Tezoatlán De Segura Y Luna"
    "B575" = "Tlacolula De Matamoros"
    "B576" = "Totontepec Villa De Morelos"
    "B577" = "Villa De Tututepec De Melchor Ocampo"
    "B578" = "Zapotitlán Del Río"
    "B579" = "Zimatlán De Álvarez"
    "B596" = "Izúcar De Matamoros"
    "B601" = "Palmar De Bravo"
    "B606" = "San Salvador El Seco"
    "B613" = "Tepexi De Rodríguez"
    "B615" = "Tlacotepec De Benito Juárez"
    "B626" = "Amealco De Bonfil"
    "B628" = "Cadereyta De Montes"
    "B631" = "Jalpan De Serra"
    "B632" = "Landa De Matamoros"
    "B634" = "Pinal De Amoles"
    "B637" = "San Juan Del Río"
    "B641" = "Armadillo De Los Infante"
    "B643" = "Ciudad Del Maíz"
    "B645" = "Mexquitic De Carmona"
    "B649" = "Soledad De Graciano Sánchez"
    "B652" = "Villa De Guadalupe"
    "B653" = "Villa De Ramos"
    "B654" = "Villa De Reyes"
    "B687" = "Nacozari De García"
    "B715" = "Ixtacuixtla De Mariano Matamoros"
    "B728" = "Amatlán De Los Reyes"
    "B739" = "Ignacio De La Llave"
    "B743" = "Juchique De Ferrer"
    "B746" = "Lerdo De Tejada"
    "B748" = "Martínez De La Torre"
    "B754" = "Paso De Ovejas"
    "B761" = "Tatahuicapan De Juárez"
    "B767" = "Vega De Alatorre"
    "B786" = "Mezquital Del Oro"
    "B790" = "Nochistlán De Mejía"
    "B791" = "Noria De Ángeles"
    "B799" = "Teúl De González Ortega"
    "B800" = "Tlaltenango De Sánchez Román"
    "B802" = "Villa De Cos"
}
foreach ($ref in $nameFixes.Keys) {
    $ws.Range($ref).Value = $nameFixes[$ref]
}

# 3. Correct floating point rounding on D662 (percentage of total)
$ws.Range("D662").Value = 0.00941446613088404

# 4. Remove trailing footnote rows (811-815), shifting dimension back to A1:D809
$ws.Range("A811:A815").EntireRow.Delete() | Out-Null

